$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. In the source data it
# belongs right after the existing row that is currently row 9 (in date
# order among the unsorted rows), so insert a fresh row there and push the
# old row 9 (and everything below it) down by one - exactly like Excel's
# "Insert" on a row heading does.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data.
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44741
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112035
$ws.Range("G9").Value = "Bruselas (repollito)"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 18800
$ws.Range("N9").Value = "$/malla 15 kilos"
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 1253
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = "Hortaliza"
